$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.1764303333333333
$ws.Range("H2").Value = 0.529291
$ws.Range("I2").Value = 0.04559680146739255
$ws.Range("J2").Value = 0.04559680146739255
$ws.Range("M2").Value = 5.125314333333333
$ws.Range("N2").Value = 15.375943
$ws.Range("O2").Value = 0.08826884604112728
$ws.Range("P2").Value = 0.0882688460411273
$ws.Range("Q2").Value = 0.904260916268111
$ws.Range("R2").Value = 8.138348246412999
$ws.Range("S2").Value = 0.00402477704869312
$ws.Range("T2").Value = 0.00402477704869312
$ws.Range("G3").Value = 0.1764303333333333
$ws.Range("H3").Value = 0.529291
$ws.Range("I3").Value = 0.04559680146739255
$ws.Range("J3").Value = 0.04559680146739255
$ws.Range("O3").Value = 0.2314334425203181
$ws.Range("P3").Value = 0.2314334425203182
$ws.Range("Q3").Value = 2.370895578390111
$ws.Range("R3").Value = 21.338060205511
$ws.Range("S3").Value = 0.01055262473151415
$ws.Range("T3").Value = 0.01055262473151415
$ws.Range("G4").Value = 0.1764303333333333
$ws.Range("H4").Value = 0.529291
$ws.Range("I4").Value = 0.04559680146739255
$ws.Range("J4").Value = 0.04559680146739255
$ws.Range("M4").Value = 37.71310533333334
$ws.Range("N4").Value = 113.139316
$ws.Range("O4").Value = 0.6495001227048286
$ws.Range("P4").Value = 0.6495001227048286
$ws.Range("Q4").Value = 6.653735744995112
$ws.Range("R4").Value = 59.883621704956
$ws.Range("S4").Value = 0.02961512814801917
$ws.Range("T4").Value = 0.02961512814801917
$ws.Range("G5").Value = 0.1764303333333333
$ws.Range("H5").Value = 0.529291
$ws.Range("I5").Value = 0.04559680146739255
$ws.Range("J5").Value = 0.04559680146739255
$ws.Range("M5").Value = 1.788256333333333
$ws.Range("N5").Value = 5.364769
$ws.Range("O5").Value = 0.03079758873372595
$ws.Range("P5").Value = 0.03079758873372595
$ws.Range("Q5").Value = 0.3155026609754444
$ws.Range("R5").Value = 2.839523948779
$ws.Range("S5").Value = 0.001404271539166107
$ws.Range("T5").Value = 0.001404271539166108
$ws.Range("I6").Value = 0.9386165989824495
$ws.Range("J6").Value = 0.9386165989824495
$ws.Range("M6").Value = 5.125314333333333
$ws.Range("N6").Value = 15.375943
$ws.Range("O6").Value = 0.08826884604112728
$ws.Range("P6").Value = 0.0882688460411273
$ws.Range("Q6").Value = 18.61433869275444
$ws.Range("R6").Value = 167.52904823479
$ws.Range("S6").Value = 0.08285060406722834
$ws.Range("T6").Value = 0.08285060406722836
$ws.Range("I7").Value = 0.9386165989824495
$ws.Range("J7").Value = 0.9386165989824495
$ws.Range("O7").Value = 0.2314334425203181
$ws.Range("P7").Value = 0.2314334425203182
$ws.Range("S7").Value = 0.2172272707092212
$ws.Range("T7").Value = 0.2172272707092213
$ws.Range("I8").Value = 0.9386165989824495
$ws.Range("J8").Value = 0.9386165989824495
$ws.Range("M8").Value = 37.71310533333334
$ws.Range("N8").Value = 113.139316
$ws.Range("O8").Value = 0.6495001227048286
$ws.Range("P8").Value = 0.6495001227048286
$ws.Range("Q8").Value = 136.9680901841645
$ws.Range("R8").Value = 1232.71281165748
$ws.Range("S8").Value = 0.6096315962118899
$ws.Range("T8").Value = 0.6096315962118899
$ws.Range("I9").Value = 0.9386165989824495
$ws.Range("J9").Value = 0.9386165989824495
$ws.Range("M9").Value = 1.788256333333333
$ws.Range("N9").Value = 5.364769
$ws.Range("O9").Value = 0.03079758873372595
$ws.Range("P9").Value = 0.03079758873372595
$ws.Range("Q9").Value = 6.494666842507778
$ws.Range("R9").Value = 58.45200158257001
$ws.Range("S9").Value = 0.02890712799411005
$ws.Range("T9").Value = 0.02890712799411006
$ws.Range("G10").Value = 0.061084
$ws.Range("H10").Value = 0.183252
$ws.Range("I10").Value = 0.01578659955015789
$ws.Range("J10").Value = 0.01578659955015789
$ws.Range("M10").Value = 5.125314333333333
$ws.Range("N10").Value = 15.375943
$ws.Range("O10").Value = 0.08826884604112728
$ws.Range("P10").Value = 0.0882688460411273
$ws.Range("Q10").Value = 0.3130747007373333
$ws.Range("R10").Value = 2.817672306636
$ws.Range("S10").Value = 0.001393464925205816
$ws.Range("T10").Value = 0.001393464925205817
$ws.Range("G11").Value = 0.061084
$ws.Range("H11").Value = 0.183252
$ws.Range("I11").Value = 0.01578659955015789
$ws.Range("J11").Value = 0.01578659955015789
$ws.Range("O11").Value = 0.2314334425203181
$ws.Range("P11").Value = 0.2314334425203182
$ws.Range("Q11").Value = 0.8208553641213332
$ws.Range("R11").Value = 7.387698277091999
$ws.Range("S11").Value = 0.003653547079582747
$ws.Range("T11").Value = 0.003653547079582748
$ws.Range("G12").Value = 0.061084
$ws.Range("H12").Value = 0.183252
$ws.Range("I12").Value = 0.01578659955015789
$ws.Range("J12").Value = 0.01578659955015789
$ws.Range("M12").Value = 37.71310533333334
$ws.Range("N12").Value = 113.139316
$ws.Range("O12").Value = 0.6495001227048286
$ws.Range("P12").Value = 0.6495001227048286
$ws.Range("Q12").Value = 2.303667326181333
$ws.Range("R12").Value = 20.733005935632
$ws.Range("S12").Value = 0.01025339834491954
$ws.Range("T12").Value = 0.01025339834491954
$ws.Range("G13").Value = 0.061084
$ws.Range("H13").Value = 0.183252
$ws.Range("I13").Value = 0.01578659955015789
$ws.Range("J13").Value = 0.01578659955015789
$ws.Range("M13").Value = 1.788256333333333
$ws.Range("N13").Value = 5.364769
$ws.Range("O13").Value = 0.03079758873372595
$ws.Range("P13").Value = 0.03079758873372595
$ws.Range("Q13").Value = 0.1092338498653333
$ws.Range("R13").Value = 0.983104648788
$ws.Range("S13").Value = 0.0004861892004497858
$ws.Range("T13").Value = 0.0004861892004497859
